$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 658352313258.7543
$ws.Range("C3").Value = 213442894232.639
$ws.Range("C4").Value = 84286479133.81744
$ws.Range("C5").Value = 34164644572.66166
$ws.Range("C6").Value = 29854210078.56976
$ws.Range("C7").Value = 25147468736.90444
$ws.Range("C8").Value = 12822792468.38154
$ws.Range("C9").Value = 9809135720.351103
$ws.Range("C10").Value = 9356571675.282593
$ws.Range("C11").Value = 8193279980.061813
$ws.Range("C12").Value = 7339386720.365101
$ws.Range("C13").Value = 5842478283.745811
$ws.Range("C14").Value = 5559553646.899556
$ws.Range("C15").Value = 5502725656.295998
$ws.Range("C16").Value = 5346120983.163179
